# Auto-generated Excel COM-interop script applying scheduled market-price updates
# to the Mateus_Profits workbook (per sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 172
$ws.Range("I4").Value = 178
$ws.Range("J4").Value = 142
$ws.Range("K4").Value = 178
$ws.Range("L4").Value = 142
$ws.Range("M4").Value = -64
$ws.Range("N4").Value = -370

$ws.Range("H41").Value = 343.5
$ws.Range("I41").Value = 329.375
$ws.Range("K41").Value = 329.375
$ws.Range("M41").Value = 110.625

$ws.Range("H62").Value = 5985.875
$ws.Range("I62").Value = 7847
$ws.Range("J62").Value = 4124.75
$ws.Range("K62").Value = 7847
$ws.Range("L62").Value = 4124.75
$ws.Range("M62").Value = -7223
$ws.Range("N62").Value = -5372.75

$ws.Range("H65").Value = 5985.875
$ws.Range("I65").Value = 7847
$ws.Range("J65").Value = 4124.75
$ws.Range("K65").Value = 39235
$ws.Range("L65").Value = 20623.75
$ws.Range("M65").Value = -36115
$ws.Range("N65").Value = -26863.75

$ws.Range("H70").Value = 3632.4707
$ws.Range("I70").Value = 2550
$ws.Range("J70").Value = 3864.4285
$ws.Range("K70").Value = 7650
$ws.Range("L70").Value = 11593.2855
$ws.Range("M70").Value = -7380
$ws.Range("N70").Value = -12133.2855

$ws.Range("H73").Value = 3632.4707
$ws.Range("I73").Value = 2550
$ws.Range("J73").Value = 3864.4285
$ws.Range("K73").Value = 7650
$ws.Range("L73").Value = 11593.2855
$ws.Range("M73").Value = -6714
$ws.Range("N73").Value = -13465.2855

$ws.Range("H76").Value = 16214.223
$ws.Range("I76").Value = 4241.75
$ws.Range("J76").Value = 17710.781
$ws.Range("K76").Value = 4241.75
$ws.Range("L76").Value = 17710.781
$ws.Range("M76").Value = -3926.75
$ws.Range("N76").Value = -18340.781

$ws.Range("H79").Value = 16214.223
$ws.Range("I79").Value = 4241.75
$ws.Range("J79").Value = 17710.781
$ws.Range("K79").Value = 4241.75
$ws.Range("L79").Value = 17710.781
$ws.Range("M79").Value = -3149.75
$ws.Range("N79").Value = -19894.781

$ws.Range("H138").Value = 2767.3242
$ws.Range("I138").Value = 2100.2
$ws.Range("J138").Value = 3552.1765
$ws.Range("K138").Value = 6300.599999999999
$ws.Range("L138").Value = 10656.5295
$ws.Range("M138").Value = -1160.599999999999
$ws.Range("N138").Value = -20936.5295

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3702.8147
$ws.Range("I74").Value = 2284.5715
$ws.Range("K74").Value = 2284.5715
$ws.Range("M74").Value = -1410.5715

$ws.Range("H77").Value = 3702.8147
$ws.Range("I77").Value = 2284.5715
$ws.Range("K77").Value = 11422.8575
$ws.Range("M77").Value = -7054.8575

$ws.Range("H102").Value = 3655.7368
$ws.Range("I102").Value = 2797.2666
$ws.Range("K102").Value = 2797.2666
$ws.Range("M102").Value = -1175.2666

$ws.Range("H122").Value = 3036.4119
$ws.Range("I122").Value = 2401.7856
$ws.Range("K122").Value = 7205.3568
$ws.Range("M122").Value = -4755.3568

$ws.Range("H132").Value = 3779.6743
$ws.Range("I132").Value = 3145.4736
$ws.Range("K132").Value = 9436.4208
$ws.Range("M132").Value = -6906.4208

$ws.Range("H134").Value = 114799.4
$ws.Range("I134").Value = 114000
$ws.Range("K134").Value = 114000
$ws.Range("M134").Value = -108930

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 6250
$ws.Range("I8").Value = 6250
$ws.Range("K8").Value = 6250
$ws.Range("M8").Value = -6110

$ws.Range("H20").Value = 3122
$ws.Range("I20").Value = 6000
$ws.Range("J20").Value = 2162.6667
$ws.Range("K20").Value = 6000
$ws.Range("L20").Value = 2162.6667
$ws.Range("M20").Value = -5753
$ws.Range("N20").Value = -2656.6667

$ws.Range("H99").Value = 3880.8333
$ws.Range("I99").Value = 2311.842
$ws.Range("J99").Value = 6590.909
$ws.Range("K99").Value = 2311.842
$ws.Range("L99").Value = 6590.909
$ws.Range("M99").Value = -813.8420000000001
$ws.Range("N99").Value = -9586.909

$ws.Range("H134").Value = 4394
$ws.Range("I134").Value = 4423.913
$ws.Range("K134").Value = 13271.739
$ws.Range("M134").Value = -10736.739

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3263.7646
$ws.Range("I31").Value = 2493.3333
$ws.Range("J31").Value = 4364.381
$ws.Range("K31").Value = 2493.3333
$ws.Range("L31").Value = 4364.381
$ws.Range("M31").Value = -2198.3333
$ws.Range("N31").Value = -4954.381

$ws.Range("H34").Value = 3263.7646
$ws.Range("I34").Value = 2493.3333
$ws.Range("J34").Value = 4364.381
$ws.Range("K34").Value = 2493.3333
$ws.Range("L34").Value = 4364.381
$ws.Range("M34").Value = -2291.3333
$ws.Range("N34").Value = -4768.381

$ws.Range("H41").Value = 10806.333
$ws.Range("I41").Value = 2861.1667
$ws.Range("J41").Value = 26696.666
$ws.Range("K41").Value = 2861.1667
$ws.Range("L41").Value = 26696.666
$ws.Range("M41").Value = -2433.1667
$ws.Range("N41").Value = -27552.666

$ws.Range("H94").Value = 1848.5
$ws.Range("J94").Value = 1831.8889
$ws.Range("L94").Value = 1831.8889
$ws.Range("N94").Value = -2733.8889

$ws.Range("H105").Value = 2107.0833
$ws.Range("I105").Value = 1755.4286
$ws.Range("K105").Value = 1755.4286
$ws.Range("M105").Value = -8.42859999999996

$ws.Range("H122").Value = 2915.9565
$ws.Range("I122").Value = 2912.5789
$ws.Range("K122").Value = 8737.736699999999
$ws.Range("M122").Value = -6287.736699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2167.0386
$ws.Range("I5").Value = 1897.0714
$ws.Range("J5").Value = 2482
$ws.Range("K5").Value = 5691.2142
$ws.Range("L5").Value = 7446
$ws.Range("M5").Value = -5579.2142
$ws.Range("N5").Value = -7670

$ws.Range("H92").Value = 102.181816
$ws.Range("I92").Value = 46.333332
$ws.Range("K92").Value = 138.999996
$ws.Range("M92").Value = 1109.000004

$ws.Range("H121").Value = 16666830
$ws.Range("I121").Value = 195.4
$ws.Range("J121").Value = 100000000
$ws.Range("K121").Value = 586.2
$ws.Range("L121").Value = 300000000
$ws.Range("M121").Value = 723.8
$ws.Range("N121").Value = -300002620

$ws.Range("H132").Value = 116666936
$ws.Range("I132").Value = 250000050
$ws.Range("J132").Value = 50000376
$ws.Range("K132").Value = 2250000450
$ws.Range("L132").Value = 450003384
$ws.Range("M132").Value = -2249997920
$ws.Range("N132").Value = -450008444

$ws.Range("H135").Value = 2167.0386
$ws.Range("I135").Value = 1897.0714
$ws.Range("J135").Value = 2482
$ws.Range("K135").Value = 17073.6426
$ws.Range("L135").Value = 22338
$ws.Range("M135").Value = -14538.6426
$ws.Range("N135").Value = -27408

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13680.353
$ws.Range("I70").Value = 10322
$ws.Range("J70").Value = 16665.555
$ws.Range("K70").Value = 10322
$ws.Range("L70").Value = 16665.555
$ws.Range("M70").Value = -10052
$ws.Range("N70").Value = -17205.555

$ws.Range("H73").Value = 13680.353
$ws.Range("I73").Value = 10322
$ws.Range("J73").Value = 16665.555
$ws.Range("K73").Value = 10322
$ws.Range("L73").Value = 16665.555
$ws.Range("M73").Value = -9386
$ws.Range("N73").Value = -18537.555

$ws.Range("H102").Value = 2684.375
$ws.Range("I102").Value = 2172.647
$ws.Range("K102").Value = 2172.647
$ws.Range("M102").Value = -550.6469999999999

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H122").Value = 4463.846
$ws.Range("I122").Value = 5233.2
$ws.Range("K122").Value = 15699.6
$ws.Range("M122").Value = -13249.6

$ws.Range("H126").Value = 4378.9165
$ws.Range("I126").Value = 4454.8
$ws.Range("J126").Value = 3999.5
$ws.Range("K126").Value = 13364.4
$ws.Range("L126").Value = 11998.5
$ws.Range("M126").Value = -10894.4
$ws.Range("N126").Value = -16938.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13151.3
$ws.Range("I7").Value = 9491.714
$ws.Range("K7").Value = 9491.714
$ws.Range("M7").Value = -9379.714

$ws.Range("H22").Value = 2514.125
$ws.Range("J22").Value = 2400
$ws.Range("L22").Value = 2400
$ws.Range("N22").Value = -2990

$ws.Range("H27").Value = 2514.125
$ws.Range("J27").Value = 2400
$ws.Range("L27").Value = 2400
$ws.Range("N27").Value = -2614

$ws.Range("H40").Value = 3013.4285
$ws.Range("I40").Value = 2184.5715
$ws.Range("K40").Value = 2184.5715
$ws.Range("M40").Value = -2048.5715

$ws.Range("H122").Value = 5006
$ws.Range("I122").Value = 4182.3335
$ws.Range("J122").Value = 5623.75
$ws.Range("K122").Value = 12547.0005
$ws.Range("L122").Value = 16871.25
$ws.Range("M122").Value = -10097.0005
$ws.Range("N122").Value = -21771.25

$ws.Range("H126").Value = 13151.3
$ws.Range("I126").Value = 9491.714
$ws.Range("K126").Value = 28475.142
$ws.Range("M126").Value = -26005.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 358332.5
$ws.Range("I4").Value = 358332.5
$ws.Range("K4").Value = 358332.5
$ws.Range("M4").Value = -358219.5

$ws.Range("H62").Value = 3167.818
$ws.Range("I62").Value = 2279.8
$ws.Range("K62").Value = 2279.8
$ws.Range("M62").Value = -1655.8

$ws.Range("H65").Value = 3167.818
$ws.Range("I65").Value = 2279.8
$ws.Range("K65").Value = 11399
$ws.Range("M65").Value = -8279

$ws.Range("H124").Value = 69696
$ws.Range("J124").Value = 69696
$ws.Range("L124").Value = 69696
$ws.Range("N124").Value = -79516
